$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 324-334: update "Patients COVID-19 aux SI total" (E) and
# "Patients COVID-19 intubés" (F) counts by +1 each; the H column
# (Total hospitalisations) recalculates automatically via its formula.
$ws.Range("E324").Value = 11
$ws.Range("F324").Value = 9

$ws.Range("E325").Value = 12
$ws.Range("F325").Value = 8

$ws.Range("E326").Value = 12
$ws.Range("F326").Value = 10

$ws.Range("E327").Value = 12
$ws.Range("F327").Value = 11

$ws.Range("E328").Value = 12
$ws.Range("F328").Value = 10

$ws.Range("E329").Value = 14
$ws.Range("F329").Value = 13

$ws.Range("E330").Value = 13
$ws.Range("F330").Value = 10

$ws.Range("E331").Value = 14
$ws.Range("F331").Value = 10

$ws.Range("E332").Value = 15
$ws.Range("F332").Value = 9

$ws.Range("E333").Value = 11
$ws.Range("F333").Value = 9

$ws.Range("E334").Value = 12
$ws.Range("F334").Value = 9

# Row 335: new daily figures (cases, SI patients, intubated, hospital
# deaths, extra-hospital deaths); B, H, J, K recalc via formulas.
$ws.Range("C335").Value = 140
$ws.Range("E335").Value = 14
$ws.Range("F335").Value = 12

# Row 336: new daily figures.
$ws.Range("C336").Value = 65
$ws.Range("E336").Value = 14
$ws.Range("F336").Value = 10

# Row 337: this day previously had no data entered yet; now fill it in.
$ws.Range("C337").Value = 12
$ws.Range("E337").Value = 16
$ws.Range("F337").Value = 11
$ws.Range("G337").Value = 130

# Columns L and M ("Nb nouveaux décès à l'hôpital" / "... extra-hospitaliers")
# are formatted as Text (@) in this sheet, so a direct numeric .Value
# assignment would be stored as text (matching real Excel's behaviour for
# Text-formatted cells). Temporarily switch those cells to a numeric format
# so the typed figures land as real numbers, then restore the original
# Text format, same as the rest of the column.
$deathCells = $ws.Range("L335:M337")
$deathCells.NumberFormat = "General"
$ws.Range("L335").Value = 2
$ws.Range("M335").Value = 1
$ws.Range("L336").Value = 1
$ws.Range("L337").Value = 1
$ws.Range("M337").Value = 0
$deathCells.NumberFormat = "@"
